# ---------------------------------------------------------------------
# RouteVisualizationData.xlsx -- add new 'Oss' delivery locations
# ---------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Target values for rows 2-27 of columns A (Start Address) and B (End
# Address). Column A is written in full first, then column B, so that
# brand-new address strings land in the shared-string table in the same
# order the source workbook uses them.
$colA = @(
    'Havikstraat 20 5348 XX, Oss',
    'Goudplevier 79 5348 ZB, Oss',
    'Goudplevier 79 5348 ZB, Oss',
    'Leeuwerikstraat 5 5348 XA, Oss',
    'Marius de Langenstraat 31 5348 AK, Oss',
    'Marius de Langenstraat 31 5348 AK, Oss',
    'Marius de Langenstraat 31 5348 AK, Oss',
    'Verdistraat 350 5343 VN, Oss',
    'Staringstraat 320 5343 GN, Oss',
    'Staringstraat 320 5343 GN, Oss',
    'Hoefstraat 11 5373 KJ, Herpen',
    'Hoefstraat 11 5373 KJ, Herpen',
    'Hoefstraat 11 5373 KJ, Herpen',
    'Hoefstraat 11 5373 KJ, Herpen',
    'Hoefstraat 11 5373 KJ, Herpen',
    'Hoefstraat 11 5373 KJ, Herpen',
    'Da Costastraat 21 5343 JN, Oss',
    'Luzacstraat 20 5344 KS, Oss',
    'Luzacstraat 20 5344 KS, Oss',
    'Slingenbergstraat 10 5344 KL, Oss',
    'Basillius van Bruggelaan 2 5363 VA, Velp-Grave',
    'Basillius van Bruggelaan 2 5363 VA, Velp-Grave',
    'Basillius van Bruggelaan 2 5363 VA, Velp-Grave',
    'Basillius van Bruggelaan 2 5363 VA, Velp-Grave',
    'Basillius van Bruggelaan 2 5363 VA, Velp-Grave',
    'Basillius van Bruggelaan 2 5363 VA, Velp-Grave'
)

$colB = @(
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss',
    'Markermeer 1 5347 JM, Oss'
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}

# The duplicate-detection conditional formatting that used to cover the
# old, now fully-replaced data range is no longer meaningful; drop it.
$ws.Cells.FormatConditions.Delete()

# Leave the cursor where the author left it after the edit.
$ws.Range("B11").Select()
